$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 350; everything from row 350 downward shifts down by one.
$ws.Rows.Item(350).Insert()

# New weekly data point for the newly inserted row 350.
$ws.Cells.Item(350, 1).Value = 7
$ws.Cells.Item(350, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(350, 3).Value = "Ñuble"
$ws.Cells.Item(350, 4).Value = 45194
$ws.Cells.Item(350, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(350, 5).Value = 16
$ws.Cells.Item(350, 6).Value = 100112043
$ws.Cells.Item(350, 7).Value = "Pepino ensalada"
$ws.Cells.Item(350, 8).Value = "Sin especificar"
$ws.Cells.Item(350, 9).Value = "Primera"
$ws.Cells.Item(350, 10).Value = 60
$ws.Cells.Item(350, 11).Value = 12000
$ws.Cells.Item(350, 12).Value = 12000
$ws.Cells.Item(350, 13).Value = 12000
$ws.Cells.Item(350, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(350, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(350, 16).Value = 200
$ws.Cells.Item(350, 17).Value = 60
$ws.Cells.Item(350, 18).Value = "Hortaliza"
